$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Aimil"
$ws.Range("B2").Value = "m"

$c2 = $ws.Range("C2")
$c2.NumberFormat = "@"
$c2.Value = "11/04/1963"
$c2.ClearFormats()

$ws.Range("D2").Value = "620 Mesta Way"
$ws.Range("E2").Value = "Scottsdale"
$ws.Range("F2").Value = "Arizona"
$ws.Range("G2").Value = 863266
$ws.Range("H2").Value = 4806825343
$ws.Range("I2").Value = "amcphatere@github.io"
$ws.Range("J2").Value = "4khIDBFT5L"
